{"js": "// Replace the multiplication-table equations in the body with the new\n// set of equations. Each \"before\" equation text is unique in the\n// document, so a simple search+replace per pair is sufficient and keeps\n// run/paragraph formatting (rFonts, sz, etc.) untouched.\nconst replacements = [\n  [\"97\u00d747=\", \"67\u00d733=\"],\n  [\"37\u00d745=\", \"44\u00d743=\"],\n  [\"43\u00d749=\", \"21\u00d784=\"],\n  [\"34\u00d748=\", \"33\u00d771=\"],\n  [\"18\u00d757=\", \"48\u00d794=\"],\n  [\"75\u00d749=\", \"52\u00d794=\"],\n  [\"20\u00d719=\", \"75\u00d750=\"],\n  [\"51\u00d770=\", \"15\u00d764=\"],\n  [\"86\u00d734=\", \"72\u00d768=\"],\n  [\"54\u00d752=\", \"45\u00d754=\"],\n  [\"56\u00d724=\", \"54\u00d724=\"],\n  [\"59\u00d787=\", \"55\u00d741=\"],\n  [\"73\u00d756=\", \"38\u00d765=\"],\n  [\"61\u00d737=\", \"55\u00d734=\"],\n  [\"52\u00d741=\", \"88\u00d788=\"],\n  [\"69\u00d726=\", \"52\u00d727=\"],\n  [\"47\u00d770=\", \"65\u00d774=\"],\n  [\"33\u00d734=\", \"85\u00d767=\"],\n  [\"78\u00d719=\", \"48\u00d725=\"],\n  [\"69\u00d769=\", \"47\u00d781=\"],\n  [\"23\u00d792=\", \"35\u00d725=\"],\n  [\"86\u00d756=\", \"90\u00d781=\"],\n  [\"14\u00d796=\", \"30\u00d717=\"],\n  [\"27\u00d751=\", \"67\u00d738=\"],\n  [\"40\u00d779=\", \"86\u00d793=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-table equations in the document body with the\n# new set of equations. Each \"before\" equation text is unique in the\n# document, so Find/Replace (wdReplaceAll = 2) per pair is sufficient and\n# leaves run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"97\u00d747=\", \"67\u00d733=\"),\n    @(\"37\u00d745=\", \"44\u00d743=\"),\n    @(\"43\u00d749=\", \"21\u00d784=\"),\n    @(\"34\u00d748=\", \"33\u00d771=\"),\n    @(\"18\u00d757=\", \"48\u00d794=\"),\n    @(\"75\u00d749=\", \"52\u00d794=\"),\n    @(\"20\u00d719=\", \"75\u00d750=\"),\n    @(\"51\u00d770=\", \"15\u00d764=\"),\n    @(\"86\u00d734=\", \"72\u00d768=\"),\n    @(\"54\u00d752=\", \"45\u00d754=\"),\n    @(\"56\u00d724=\", \"54\u00d724=\"),\n    @(\"59\u00d787=\", \"55\u00d741=\"),\n    @(\"73\u00d756=\", \"38\u00d765=\"),\n    @(\"61\u00d737=\", \"55\u00d734=\"),\n    @(\"52\u00d741=\", \"88\u00d788=\"),\n    @(\"69\u00d726=\", \"52\u00d727=\"),\n    @(\"47\u00d770=\", \"65\u00d774=\"),\n    @(\"33\u00d734=\", \"85\u00d767=\"),\n    @(\"78\u00d719=\", \"48\u00d725=\"),\n    @(\"69\u00d769=\", \"47\u00d781=\"),\n    @(\"23\u00d792=\", \"35\u00d725=\"),\n    @(\"86\u00d756=\", \"90\u00d781=\"),\n    @(\"14\u00d796=\", \"30\u00d717=\"),\n    @(\"27\u00d751=\", \"67\u00d738=\"),\n    @(\"40\u00d779=\", \"86\u00d793=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
